$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DirectoryResults")

$ws.Range("B1").Value = "directory_result"
$ws.Range("B2").Value = " Thomas Abowd"
$ws.Range("B3").Value = " Danielle Abrams"
$ws.Range("B4").Value = " Dany Adams, Ph.D."
$ws.Range("B5").Value = " Kristina Aikens, PhD"
